$d = $word.ActiveDocument
$d.Content.Find.Execute("wep", $true, $true, $false, $false, $false,
                         $true, 1, $false, "web", 2)
